# Edit: "Added plugins filtering by prefix and plugin development."
#
# This script updates the first worksheet ("Semanticke clenenie") of the
# MQTT_topics workbook:
#   - the "system" device category is renamed to "server" and gets two
#     extra rows (a "perc" measure row and a "cmd"/"period" row plus the
#     usual "---"/"---" STATUS,RESET command row)
#   - the "iot" device category is renamed to "iotgate" and gets an extra
#     "cmd"/"period" row
#   - the long description of the "category" column is reworded
#   - the sheet grows from 21 to 24 used rows and the selection moves to
#     the newly added comment block

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- grow the data table: insert 3 new rows inside the existing block ---
# (row 15/16 are already blank in the source sheet, so only 3 explicit
# inserts are required to reach the final 17-row data block)
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(6).Insert()

# --- row 1 : headers (unchanged, left as-is) ---

# --- data rows 2-18 : full target content -----------------------------
# row 2 : server / data / temp / val / 60
$ws.Cells.Item(2,1).Value = "server"
$ws.Cells.Item(2,2).Value = "data"
$ws.Cells.Item(2,3).Value = "temp"
$ws.Cells.Item(2,4).Value = "val"
$ws.Cells.Item(2,5).Value = 60

# row 3 : server / data / temp / perc / 80   (new row)
$ws.Cells.Item(3,1).Value = "server"
$ws.Cells.Item(3,2).Value = "data"
$ws.Cells.Item(3,3).Value = "temp"
$ws.Cells.Item(3,4).Value = "perc"
$ws.Cells.Item(3,5).Value = 80

# row 4 : server / state / temp / max / 75
$ws.Cells.Item(4,1).Value = "server"
$ws.Cells.Item(4,2).Value = "state"
$ws.Cells.Item(4,3).Value = "temp"
$ws.Cells.Item(4,4).Value = "max"
$ws.Cells.Item(4,5).Value = 75

# row 5 : server / cmd / period / val / 5   (new row)
$ws.Cells.Item(5,1).Value = "server"
$ws.Cells.Item(5,2).Value = "cmd"
$ws.Cells.Item(5,3).Value = "period"
$ws.Cells.Item(5,4).Value = "val"
$ws.Cells.Item(5,5).Value = 5

# row 6 : server / cmd / --- / --- / " STATUS, RESET"   (new row)
$ws.Cells.Item(6,1).Value = "server"
$ws.Cells.Item(6,2).Value = "cmd"
$ws.Cells.Item(6,3).Value = "---"
$ws.Cells.Item(6,4).Value = "---"
$ws.Cells.Item(6,5).Value = " STATUS, RESET"

# row 7 : fan / state / --- / --- / "Online, Offline, Active, Idle"
$ws.Cells.Item(7,1).Value = "fan"
$ws.Cells.Item(7,2).Value = "state"
$ws.Cells.Item(7,3).Value = "---"
$ws.Cells.Item(7,4).Value = "---"
$ws.Cells.Item(7,5).Value = "Online, Offline, Active, Idle"

# row 8 : fan / state / percon / / 95
$ws.Cells.Item(8,1).Value = "fan"
$ws.Cells.Item(8,2).Value = "state"
$ws.Cells.Item(8,3).Value = "percon"
$ws.Cells.Item(8,5).Value = 95

# row 9 : fan / state / percoff / / 60
$ws.Cells.Item(9,1).Value = "fan"
$ws.Cells.Item(9,2).Value = "state"
$ws.Cells.Item(9,3).Value = "percoff"
$ws.Cells.Item(9,5).Value = 60

# row 10 : fan / state / tempon / / 71.25
$ws.Cells.Item(10,1).Value = "fan"
$ws.Cells.Item(10,2).Value = "state"
$ws.Cells.Item(10,3).Value = "tempon"
$ws.Cells.Item(10,5).Value = 71.25

# row 11 : fan / state / tempoff / / 45
$ws.Cells.Item(11,1).Value = "fan"
$ws.Cells.Item(11,2).Value = "state"
$ws.Cells.Item(11,3).Value = "tempoff"
$ws.Cells.Item(11,5).Value = 45

# row 12 : fan / cmd / --- / --- / "ON, OFF, TOGGLE, STATUS, RESET"
$ws.Cells.Item(12,1).Value = "fan"
$ws.Cells.Item(12,2).Value = "cmd"
$ws.Cells.Item(12,3).Value = "---"
$ws.Cells.Item(12,4).Value = "---"
$ws.Cells.Item(12,5).Value = "ON, OFF, TOGGLE, STATUS, RESET"

# row 13 : fan / cmd / percon / / 85
$ws.Cells.Item(13,1).Value = "fan"
$ws.Cells.Item(13,2).Value = "cmd"
$ws.Cells.Item(13,3).Value = "percon"
$ws.Cells.Item(13,5).Value = 85

# row 14 : fan / cmd / percoff / / 75
$ws.Cells.Item(14,1).Value = "fan"
$ws.Cells.Item(14,2).Value = "cmd"
$ws.Cells.Item(14,3).Value = "percoff"
$ws.Cells.Item(14,5).Value = 75

# row 15 : fan / cmd / tempon / / 71.25
$ws.Cells.Item(15,1).Value = "fan"
$ws.Cells.Item(15,2).Value = "cmd"
$ws.Cells.Item(15,3).Value = "tempon"
$ws.Cells.Item(15,5).Value = 71.25

# row 16 : fan / cmd / tempoff / / 45
$ws.Cells.Item(16,1).Value = "fan"
$ws.Cells.Item(16,2).Value = "cmd"
$ws.Cells.Item(16,3).Value = "tempoff"
$ws.Cells.Item(16,5).Value = 45

# row 17 : iotgate / state / --- / --- / "Online, Offline"
$ws.Cells.Item(17,1).Value = "iotgate"
$ws.Cells.Item(17,2).Value = "state"
$ws.Cells.Item(17,3).Value = "---"
$ws.Cells.Item(17,4).Value = "---"
$ws.Cells.Item(17,5).Value = "Online, Offline"

# row 18 : iotgate / cmd / period / val / 60   (brand new row, was blank)
$ws.Cells.Item(18,1).Value = "iotgate"
$ws.Cells.Item(18,2).Value = "cmd"
$ws.Cells.Item(18,3).Value = "period"
$ws.Cells.Item(18,4).Value = "val"
$ws.Cells.Item(18,5).Value = 60

# --- reword the "category" comment block (now rows 23/24) -------------
$ws.Cells.Item(24,1).Value = "Toto sú pevné vymenované kategórie, ktoré brána rozlišuje pre každy plugin rovnako. Zodpovedajú funkcionalitám pluginov.`r`nPoznajú ich len brána a hardvérový komponent. Brána len volá im zodpovedajúce metódy pluginu."

# --- update the selection to match the freshly edited comment block ---
$ws.Range("A24:J24").Select()
